$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Fix the ROOMS006 (row 8) result from "fail" to "pass"
$ws.Range("D8").Value = "pass"

# 2. Update the signoff date for every existing data row (3-32) from 4/21/2013 to 4/23/2013
$ws.Range("C3:C32").Value = 41387

# 3. Append new row 33 for rooms031 - copy formatting from row 31 (a full A:D data row)
#    so number formats / alignment styles match, then overwrite the values.
$ws.Range("A31:D31").Copy($ws.Range("A33:D33"))
$ws.Range("A33").Value = "rooms031"
$ws.Range("B33").Value = "Jared Cox"
$ws.Range("C33").Value = 41387
$ws.Range("D33").Value = "pass"

# 4. Add the remark explaining why the file was moved to "No Error" files
$ws.Range("E33").Value = "Moved this file to No Error files.  Program overwites duplicates."
$ws.Range("E33").HorizontalAlignment = -4108
$ws.Range("E33").WrapText = $true
$ws.Range("E33").EntireRow.RowHeight = 30

# 5. Update the view state to reflect scrolling down to the newly signed-off row
[void]$ws.Range("E33").Select()
